$p = $ppt.ActivePresentation
$s = $p.Slides.Item(3)
$sh = $s.Shapes.Item(3)

$tf = $sh.TextFrame
$tr = $tf.TextRange
$para1 = $tr.Paragraphs(1, 1)

# Split "Stars" into two runs: strike-through "Stars" + new trailing note.
$null = $para1.InsertAfter(" (these can be masked out)")
$starsRange = $para1.Characters(1, 5)
$starsRange.Font.StrikeThrough = $true

# The textbox auto-fits to the new (wider) text; update its position/size to match.
$sh.Left = 202.1426
$sh.Width = 284.0636
